# Manage Interviewers flow - Done
# Updates the "AMS" sheet's sprint-run history table:
#  - Row 10: refine the already-present run (style touch-up + more precise timestamp)
#  - Row 11: add the "live_145_hotfix" run
#  - Row 12: add the "live_145_hf2" run

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMS")

# --- Row 10 : existing run, timestamp refined + formatting normalised ---
$ws.Cells.Item(10,1).Value = "2021-06-11"
$ws.Cells.Item(10,2).Value = 44358.66002978009
$ws.Cells.Item(10,3).Value = "pavan_demo_145"
$ws.Cells.Item(10,4).Value = 119
$ws.Cells.Item(10,5).Value = 117
$ws.Cells.Item(10,6).Value = 2
$ws.Cells.Item(10,7).Value = 2.95

# --- Row 11 : new run "live_145_hotfix" ---
$ws.Cells.Item(11,1).Value = "2021-06-16"
$ws.Cells.Item(11,2).Value = 44363.69445333333
$ws.Cells.Item(11,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11,3).Value = "live_145_hotfix"
$ws.Cells.Item(11,4).Value = 119
$ws.Cells.Item(11,5).Value = 111
$ws.Cells.Item(11,6).Value = 8
$ws.Cells.Item(11,7).Value = 4.24

# --- Row 12 : new run "live_145_hf2" ---
$ws.Cells.Item(12,1).Value = "2021-06-16"
$ws.Cells.Item(12,2).Value = 44363.81532318905
$ws.Cells.Item(12,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12,3).Value = "live_145_hf2"
$ws.Cells.Item(12,4).Value = 119
$ws.Cells.Item(12,5).Value = 117
$ws.Cells.Item(12,6).Value = 2
$ws.Cells.Item(12,7).Value = 2.87
